$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6630574638774661
$ws.Range("J2").Value = 0.663057463877466
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.022792
$ws.Range("N2").Value = 0.06837600000000001
$ws.Range("O2").Value = 0.001916327914826657
$ws.Range("P2").Value = 0.001916327914826657
$ws.Range("Q2").Value = 0.003715430282666667
$ws.Range("R2").Value = 0.03343887254400001
$ws.Range("S2").Value = 0.001270635527162556
$ws.Range("T2").Value = 0.001270635527162556

# Row 3
$ws.Range("I3").Value = 0.6630574638774661
$ws.Range("J3").Value = 0.663057463877466
$ws.Range("O3").Value = 0.3701235913233977
$ws.Range("P3").Value = 0.3701235913233977
$ws.Range("S3").Value = 0.2454132097841118
$ws.Range("T3").Value = 0.2454132097841118

# Row 4
$ws.Range("I4").Value = 0.6630574638774661
$ws.Range("J4").Value = 0.663057463877466
$ws.Range("M4").Value = 7.468693666666667
$ws.Range("N4").Value = 22.406081
$ws.Range("O4").Value = 0.6279600807617757
$ws.Range("P4").Value = 0.6279600807617757
$ws.Range("Q4").Value = 1.217506608507111
$ws.Range("R4").Value = 10.957559476564
$ws.Range("S4").Value = 0.4163736185661918
$ws.Range("T4").Value = 0.4163736185661917

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08283833333333333
$ws.Range("H5").Value = 0.248515
$ws.Range("I5").Value = 0.3369425361225339
$ws.Range("J5").Value = 0.3369425361225339
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.022792
$ws.Range("N5").Value = 0.06837600000000001
$ws.Range("O5").Value = 0.001916327914826657
$ws.Range("P5").Value = 0.001916327914826657
$ws.Range("Q5").Value = 0.001888051293333333
$ws.Range("R5").Value = 0.01699246164
$ws.Range("S5").Value = 0.000645692387664101
$ws.Range("T5").Value = 0.0006456923876641009

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08283833333333333
$ws.Range("H6").Value = 0.248515
$ws.Range("I6").Value = 0.3369425361225339
$ws.Range("J6").Value = 0.3369425361225339
$ws.Range("O6").Value = 0.3701235913233977
$ws.Range("P6").Value = 0.3701235913233977
$ws.Range("Q6").Value = 0.3646621853622222
$ws.Range("R6").Value = 3.28195966826
$ws.Range("S6").Value = 0.1247103815392859
$ws.Range("T6").Value = 0.1247103815392859

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08283833333333333
$ws.Range("H7").Value = 0.248515
$ws.Range("I7").Value = 0.3369425361225339
$ws.Range("J7").Value = 0.3369425361225339
$ws.Range("M7").Value = 7.468693666666667
$ws.Range("N7").Value = 22.406081
$ws.Range("O7").Value = 0.6279600807617757
$ws.Range("P7").Value = 0.6279600807617757
$ws.Range("Q7").Value = 0.6186941355238889
$ws.Range("R7").Value = 5.568247219714999
$ws.Range("S7").Value = 0.2115864621955839
$ws.Range("T7").Value = 0.2115864621955839
